$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: H3 (Absent) was 0, now Absent -> 1
$ws.Range("H3").Value = 1

# Row 5: H5 was an empty/inline-string cell, now Absent -> 0 (numeric)
$ws.Range("H5").Value = 0

# Row 12: H12 (Absent) was 0, now Absent -> 1
$ws.Range("H12").Value = 1

# Row 13: H13 was an empty/inline-string cell, now Absent -> 0 (numeric)
$ws.Range("H13").Value = 0
